# fix(gui) step 1 and 2
# INGLETES price list: bump the "updated" date stamp in A1 by one day,
# and correct the prices for the PINO (row 30 / step 1) and ALGARROBO
# (row 31 / step 2) inglete items.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1 holds the date serial (formatted as a date); move it from 45308 (2024-01-17)
# to 45309 (2024-01-18).
$ws.Range("A1").Value = 45309

# Step 1: Inglete de PINO (row 30) price update.
$ws.Range("D30").Value = 760

# Step 2: Inglete de ALGARROBO (row 31) price update.
$ws.Range("D31").Value = 1520
